# Insert a new weekly record before the existing row 108, shifting all
# subsequent rows (108..214) down by one (to 109..215), and populate the
# new row 108 with the new data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 108 (pushes rows 108-214 down to 109-215)
$ws.Rows.Item(108).Insert()

# Fill in the new row 108 with the new data
$ws.Cells.Item(108,1).Value  = 2
$ws.Cells.Item(108,2).Value  = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(108,3).Value  = "Coquimbo"
$ws.Cells.Item(108,4).Value  = 44874
$ws.Cells.Item(108,5).Value  = 4
$ws.Cells.Item(108,6).Value  = 100112031
$ws.Cells.Item(108,7).Value  = "Poroto verde"
$ws.Cells.Item(108,8).Value  = "Magnum"
$ws.Cells.Item(108,9).Value  = "Primera"
$ws.Cells.Item(108,10).Value = 500
$ws.Cells.Item(108,11).Value = 41000
$ws.Cells.Item(108,12).Value = 43000
$ws.Cells.Item(108,13).Value = 42000
$ws.Cells.Item(108,14).Value = "`$/caja 25 kilos"
$ws.Cells.Item(108,15).Value = "Provincia de Limarí"
$ws.Cells.Item(108,16).Value = 1680
$ws.Cells.Item(108,17).Value = 25
$ws.Cells.Item(108,18).Value = "Hortaliza"
